# "Generate Report for Archive" — refresh the localization status report:
# the two outstanding files have moved from "Ready for handoff" into
# "In Translation", and the Status column narrows to fit the shorter text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: Status is mirrored per-locale in columns E (zh-cn) and F (de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = $newStatus
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

# --- zh-cn sheet: Status lives in column C ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = $newStatus
$zhcn.Range("C1").ColumnWidth = 12.5

# --- de-de sheet: Status lives in column C ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = $newStatus
$dede.Range("C1").ColumnWidth = 12.5
